$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# The "Challenge Card" data table (A:cardName, B:maxCount, C:effect)
# is being reworked for the new challenge-card rules: several old
# cards (一线天/大房间/迷雾/裂隙/手层/冲突层) are removed, the
# remaining cards (大墓地/斗技场/展厅/沙虫坑) get new multi-tier
# "DL2/DL3/BOSS" effect text, and two brand-new cards are added
# (手指层/迷雾森林).
# ---------------------------------------------------------------------

# The table originally had 10 data rows (rows 2-11); the new table only
# needs 6 data rows (rows 2-7), so drop the last 4 rows entirely.
$ws.Range("A8:A11").EntireRow.Delete()

$rows = @(
    @{ Name = "大墓地"; Height = 57; Effect = "DL2：房间尺寸加1。<br>`nDL3：房间尺寸加1。墓地中的牌不受玩家控制的牌的效果影响。<br>`nBOSS：从额外牌堆将1张《命匣》牌洗入主牌堆。每次重整后，从额外牌堆将《巫妖》牌放在房间区任意槽位中。" },
    @{ Name = "斗技场"; Height = 57; Effect = "DL2：房间尺寸加1。<br>`nDL3：房间尺寸加1。战斗中，战场外的牌无法移入战场，战场中的牌无法移出战场。<br>`nBOSS ：从额外牌堆将1张《冠军》牌放在房间区任意槽位中。" },
    @{ Name = "展厅"; Height = 42.75; Effect = "DL2：房间尺寸加1。<br>`nDL3：房间尺寸加1。遗物牌使用后横置，且不会在回合结束时被复位。<br>`nBOSS：从额外牌堆将1张《收藏家》牌洗入主牌堆。" },
    @{ Name = "沙虫坑"; Height = 42.75; Effect = "DL2：房间尺寸加1。<br>`nDL3：房间尺寸加1。场上的牌无法移动到《沙虫》所在槽位中。<br>`nBOSS：从额外牌堆将3张《沙虫》牌依次放在房间区任意槽位中。" },
    @{ Name = "手指层"; Height = 42.75; Effect = "DL2：房间尺寸加1。<br>`nDL3：房间尺寸加1。玩家的手牌只能被打出。<br>`nBOSS：从额外牌堆将5张《封印手指》牌洗入主牌堆。" },
    @{ Name = "迷雾森林"; Height = 42.75; Effect = "DL2：房间尺寸加1。<br>`nDL3：房间尺寸加1。重整时，房间区中央的牌保持背面向上发出。<br>`nBOSS：" }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Name
    $c = $ws.Cells.Item($r, 3)
    $c.Value = $row.Effect
    # Make sure the effect cell uses the wrap-text style (s="1"), same as
    # the rest of the table - row 2's effect cell did not have it before.
    $c.WrapText = $true
    $ws.Rows.Item($r).RowHeight = $row.Height
    $r = $r + 1
}

# Update the visible selection to match the final sheet state.
$ws.Range("C5").Select()
